# ROADMAP.docx -> Russian translation pass
# Applies the textual changes from the target diff by locating each
# English source run via Find.Execute, then using Range.InsertXML to
# replace that run's content with the translated text while preserving
# the exact original run/paragraph XML (rPr, rsid attributes, xml:space).

$d = $word.ActiveDocument

function Replace-RunXml($searchText, $runOpenTag, $rprXml, $newText) {
    $full = $d.Content
    $found = $full.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $target = $d.Range($full.Start, $full.End)
    $runXml = $runOpenTag + $rprXml + '<w:t xml:space="preserve">' + $newText + '</w:t></w:r>'
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
           $runXml + `
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

$nbsp = [char]0x00A0
$longOld = "SmartRewards, SmartHive voting, SmartHive, Smart Web wallet , Mobile wallets, SmartNodes, Electrum, Pay to email, InstantPay,${nbsp}Core upgrades for faster syncing,${nbsp}vault, multisig (complete), trezor (hardware wallet), upgraded explorer, 24hr support center, smartnode starting from electrum, electrum smartvoting, SmartCard and Merchant Reader, Collateral change to 100k to enable better quality SmartNodes, SmartShift, SmartRewards tab in Node and Electrum wallets, SmartNode starting with Trezor, SmartCard and Merchant Reader"
$longNew = 'SmartRewards, SmartHive голосование, SmartHive, Smart Web Кошелёк, Мобильные кошельки, SmartNodes, Кошелёк Electrum, Отправка по Email, InstantPay, Обновление кода для быстрой синхронизации, Хранилище Vault, Мультиподписи, Trezor (аппаратный кошелёк), Обновлённый обозреватель блоков, Поддержка 24/7, Запуск SmartNode с кошелька Electrum, Голосование с кошелька Electrum, SmartCard и платёжное решение, Изменение залоговой суммы до 100к для улучшения качества SmartNodes, SmartShift, вкладка SmartRewards в Node и Electrum кошельках, запуск SmartNode с Trezor'
$longRunOpen = '<w:r>'
$longRpr = '<w:rPr><w:rFonts w:ascii="open sans" w:hAnsi="open sans" w:cs="open sans" w:eastAsia="Times New Roman"/><w:b/><w:bCs/><w:color w:val="2B2B2B"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
Replace-RunXml $longOld $longRunOpen $longRpr $longNew

Replace-RunXml 'Community Projects' '<w:r>' '<w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="F4B517"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' 'Проекты сообщества'
Replace-RunXml 'Projects Funded' '<w:r w:rsidRPr="00D70229">' '<w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="Times New Roman" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr>' 'Профинансировано проектов'
Replace-RunXml 'Completed Projects' '<w:r>' '<w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="F4B517"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' 'Завершённые проекты'
Replace-RunXml 'Projects Complete' '<w:r w:rsidRPr="00D70229">' '<w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="Times New Roman" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr>' 'Завершено проектов'
Replace-RunXml 'Funds Allocated' '<w:r>' '<w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/><w:color w:val="F4B517"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>' 'Выделенные средства'
Replace-RunXml 'SmartCash Invested' '<w:r w:rsidRPr="00D70229">' '<w:rPr><w:rFonts w:ascii="Tahoma" w:eastAsia="Times New Roman" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:color w:val="252525"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/></w:rPr>' 'Инвестировано средств'

Write-Output "translations applied"
